$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 62
$ws_ALC.Range("H62").Value = 7567.857
$ws_ALC.Range("I62").Value = 7162.5
$ws_ALC.Range("K62").Value = 7162.5
$ws_ALC.Range("M62").Value = -6538.5

# ALC row 65
$ws_ALC.Range("H65").Value = 7567.857
$ws_ALC.Range("I65").Value = 7162.5
$ws_ALC.Range("K65").Value = 35812.5
$ws_ALC.Range("M65").Value = -32692.5

# ALC row 86
$ws_ALC.Range("H86").Value = 173162850
$ws_ALC.Range("J86").Value = 156382240
$ws_ALC.Range("L86").Value = 156382240
$ws_ALC.Range("N86").Value = -156384486

# ALC row 89
$ws_ALC.Range("H89").Value = 173162850
$ws_ALC.Range("J89").Value = 156382240
$ws_ALC.Range("L89").Value = 781911200
$ws_ALC.Range("N89").Value = -781922432

# ARM row 32
$ws_ARM.Range("H32").Value = 18895228
$ws_ARM.Range("I32").Value = 21796826
$ws_ARM.Range("K32").Value = 21796826
$ws_ARM.Range("M32").Value = -21796539

# ARM row 45
$ws_ARM.Range("H45").Value = 2956.4
$ws_ARM.Range("I45").Value = 1812.8
$ws_ARM.Range("K45").Value = 1812.8
$ws_ARM.Range("M45").Value = -1435.8

# ARM row 46
$ws_ARM.Range("H46").Value = 3450.6667
$ws_ARM.Range("J46").Value = 3450.6667
$ws_ARM.Range("L46").Value = 3450.6667
$ws_ARM.Range("N46").Value = -4088.6667

# ARM row 61
$ws_ARM.Range("H61").Value = 4221.4644
$ws_ARM.Range("I61").Value = 4973.5386
$ws_ARM.Range("K61").Value = 4973.5386
$ws_ARM.Range("M61").Value = -4761.5386

# ARM row 74
$ws_ARM.Range("H74").Value = 2324.3137
$ws_ARM.Range("I74").Value = 1926.4043
$ws_ARM.Range("K74").Value = 1926.4043
$ws_ARM.Range("M74").Value = -1052.4043

# ARM row 77
$ws_ARM.Range("H77").Value = 2324.3137
$ws_ARM.Range("I77").Value = 1926.4043
$ws_ARM.Range("K77").Value = 9632.021499999999
$ws_ARM.Range("M77").Value = -5264.021499999999

# ARM row 122
$ws_ARM.Range("H122").Value = 5498.5
$ws_ARM.Range("I122").Value = 4000
$ws_ARM.Range("J122").Value = 6997
$ws_ARM.Range("K122").Value = 12000
$ws_ARM.Range("L122").Value = 20991
$ws_ARM.Range("M122").Value = -9550
$ws_ARM.Range("N122").Value = -25891

# ARM row 128
$ws_ARM.Range("H128").Value = 115999
$ws_ARM.Range("J128").Value = 115999
$ws_ARM.Range("L128").Value = 115999
$ws_ARM.Range("N128").Value = -125959

# ARM row 132
$ws_ARM.Range("H132").Value = 719162.2
$ws_ARM.Range("I132").Value = 1433894.1
$ws_ARM.Range("K132").Value = 4301682.300000001
$ws_ARM.Range("M132").Value = -4299152.300000001

# ARM row 136
$ws_ARM.Range("H136").Value = 4221.4644
$ws_ARM.Range("I136").Value = 4973.5386
$ws_ARM.Range("K136").Value = 14920.6158
$ws_ARM.Range("M136").Value = -12370.6158

# BSM row 11
$ws_BSM.Range("H11").Value = 233.75
$ws_BSM.Range("I11").Value = 233.75
$ws_BSM.Range("J11").Value = 0
$ws_BSM.Range("K11").Value = 233.75
$ws_BSM.Range("L11").Value = 0
$ws_BSM.Range("M11").Value = -93.75
$ws_BSM.Range("N11").ClearContents()

# BSM row 105
$ws_BSM.Range("H105").Value = 2655.2104
$ws_BSM.Range("I105").Value = 2200.6365
$ws_BSM.Range("K105").Value = 2200.6365
$ws_BSM.Range("M105").Value = -453.6365000000001

# BSM row 134
$ws_BSM.Range("H134").Value = 3972171.2
$ws_BSM.Range("I134").Value = 5496475.5
$ws_BSM.Range("K134").Value = 16489426.5
$ws_BSM.Range("M134").Value = -16486891.5

# CRP row 3
$ws_CRP.Range("H3").Value = 3002
$ws_CRP.Range("I3").Value = 3002
$ws_CRP.Range("J3").Value = 0
$ws_CRP.Range("K3").Value = 3002
$ws_CRP.Range("L3").Value = 0
$ws_CRP.Range("M3").Value = -2889
$ws_CRP.Range("N3").ClearContents()

# CRP row 16
$ws_CRP.Range("H16").Value = 675.5
$ws_CRP.Range("I16").Value = 675.5
$ws_CRP.Range("K16").Value = 675.5
$ws_CRP.Range("M16").Value = -388.5

# CRP row 68
$ws_CRP.Range("H68").Value = 83999.5
$ws_CRP.Range("I68").Value = 32000
$ws_CRP.Range("J68").Value = 135999
$ws_CRP.Range("K68").Value = 32000
$ws_CRP.Range("L68").Value = 135999
$ws_CRP.Range("M68").Value = -31251
$ws_CRP.Range("N68").Value = -137497

# CRP row 71
$ws_CRP.Range("H71").Value = 83999.5
$ws_CRP.Range("I71").Value = 32000
$ws_CRP.Range("J71").Value = 135999
$ws_CRP.Range("K71").Value = 96000
$ws_CRP.Range("L71").Value = 407997
$ws_CRP.Range("M71").Value = -92256
$ws_CRP.Range("N71").Value = -415485

# CRP row 86
$ws_CRP.Range("H86").Value = 36794.207
$ws_CRP.Range("J86").Value = 36334
$ws_CRP.Range("L86").Value = 36334
$ws_CRP.Range("N86").Value = -38580

# CRP row 89
$ws_CRP.Range("H89").Value = 36794.207
$ws_CRP.Range("J89").Value = 36334
$ws_CRP.Range("L89").Value = 181670
$ws_CRP.Range("N89").Value = -192902

# CRP row 107
$ws_CRP.Range("H107").Value = 1244.8889
$ws_CRP.Range("I107").Value = 709.6
$ws_CRP.Range("J107").Value = 1559.7646
$ws_CRP.Range("K107").Value = 709.6
$ws_CRP.Range("L107").Value = 1559.7646
$ws_CRP.Range("M107").Value = 1210.4
$ws_CRP.Range("N107").Value = -5399.7646

# CRP row 113
$ws_CRP.Range("H113").Value = 675.5
$ws_CRP.Range("I113").Value = 675.5
$ws_CRP.Range("K113").Value = 675.5
$ws_CRP.Range("M113").Value = 1494.5

# CRP row 132
$ws_CRP.Range("H132").Value = 3315.2222
$ws_CRP.Range("I132").Value = 2922.5518
$ws_CRP.Range("J132").Value = 4942
$ws_CRP.Range("K132").Value = 8767.6554
$ws_CRP.Range("L132").Value = 14826
$ws_CRP.Range("M132").Value = -6237.6554
$ws_CRP.Range("N132").Value = -19886

# GSM row 80
$ws_GSM.Range("H80").Value = 4332.3335
$ws_GSM.Range("I80").Value = 5000
$ws_GSM.Range("J80").Value = 3998.5
$ws_GSM.Range("K80").Value = 5000
$ws_GSM.Range("L80").Value = 3998.5
$ws_GSM.Range("M80").Value = -4002
$ws_GSM.Range("N80").Value = -5994.5

# GSM row 83
$ws_GSM.Range("H83").Value = 4332.3335
$ws_GSM.Range("I83").Value = 5000
$ws_GSM.Range("J83").Value = 3998.5
$ws_GSM.Range("K83").Value = 25000
$ws_GSM.Range("L83").Value = 19992.5
$ws_GSM.Range("M83").Value = -20008
$ws_GSM.Range("N83").Value = -29976.5

# GSM row 113
$ws_GSM.Range("H113").Value = 11476.9
$ws_GSM.Range("I113").Value = 1641.1111
$ws_GSM.Range("J113").Value = 99999
$ws_GSM.Range("K113").Value = 1641.1111
$ws_GSM.Range("L113").Value = 99999
$ws_GSM.Range("M113").Value = 528.8888999999999
$ws_GSM.Range("N113").Value = -104339

# GSM row 122
$ws_GSM.Range("H122").Value = 999.5
$ws_GSM.Range("I122").Value = 999.5
$ws_GSM.Range("K122").Value = 2998.5
$ws_GSM.Range("M122").Value = -548.5

# GSM row 126
$ws_GSM.Range("H126").Value = 2689.9167
$ws_GSM.Range("I126").Value = 2459.4285
$ws_GSM.Range("K126").Value = 7378.2855
$ws_GSM.Range("M126").Value = -4908.2855

# GSM row 132
$ws_GSM.Range("H132").Value = 10649.667
$ws_GSM.Range("I132").Value = 14299.667
$ws_GSM.Range("J132").Value = 6999.6665
$ws_GSM.Range("K132").Value = 42899.001
$ws_GSM.Range("L132").Value = 20998.9995
$ws_GSM.Range("M132").Value = -40369.001
$ws_GSM.Range("N132").Value = -26058.9995

# LTW row 122
$ws_LTW.Range("H122").Value = 13360.77
$ws_LTW.Range("I122").Value = 11974.25
$ws_LTW.Range("K122").Value = 35922.75
$ws_LTW.Range("M122").Value = -33472.75

# LTW row 132
$ws_LTW.Range("H132").Value = 459112.1
$ws_LTW.Range("I132").Value = 592109.8
$ws_LTW.Range("K132").Value = 1776329.4
$ws_LTW.Range("M132").Value = -1773799.4

# WVR row 3
$ws_WVR.Range("H3").Value = 626300.75
$ws_WVR.Range("I3").Value = 834401
$ws_WVR.Range("K3").Value = 834401
$ws_WVR.Range("M3").Value = -834287

# WVR row 18
$ws_WVR.Range("H18").Value = 17753.375
$ws_WVR.Range("I18").Value = 17751.5
$ws_WVR.Range("J18").Value = 17755.25
$ws_WVR.Range("K18").Value = 17751.5
$ws_WVR.Range("L18").Value = 17755.25
$ws_WVR.Range("M18").Value = -17578.5
$ws_WVR.Range("N18").Value = -18101.25

# WVR row 132
$ws_WVR.Range("H132").Value = 26435.023
$ws_WVR.Range("I132").Value = 29414.75
$ws_WVR.Range("J132").Value = 4981
$ws_WVR.Range("K132").Value = 88244.25
$ws_WVR.Range("L132").Value = 14943
$ws_WVR.Range("M132").Value = -85714.25
$ws_WVR.Range("N132").Value = -20003
